$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Insert 34 new rows above the existing data (row 2), pushing old rows 2-21 down to 36-55
$ws.Rows("2:35").Insert()
$ws.Range("A2:F35").ClearFormats()

# Force text storage (matches scraped numeric-looking string columns A, C, D, E, F)
$ws.Range("A2:A35").NumberFormat = "@"
$ws.Range("C2:D35").NumberFormat = "@"
$ws.Range("E2:F35").NumberFormat = "@"

$newData = @(
    @("3895", "1", "7", "1", "57.80%", "YES"),
    @("3897", "1", "4", "0", "25.58%", "NO"),
    @("3899", "1", "4", "2", "50.00%", "YES"),
    @("3974", "1", "1", "0", "2.25%", "NO"),
    @("3976", "", "", "", "", "NO"),
    @("3978", "", "", "", "", "NO"),
    @("4061", "4", "", "", "", "NO"),
    @("4062", "", "", "", "", "NO"),
    @("4063", "", "", "", "", "NO"),
    @("4064", "5", "0", "0", "1.87%", "NO"),
    @("4171", "4", "0", "0", "3.35%", "NO"),
    @("4173", "4", "0", "0", "", "NO"),
    @("4203", "", "", "", "", "NO"),
    @("4268", "3", "1", "0", "7.26%", "NO"),
    @("4310", "", "", "", "", "NO"),
    @("4316", "6", "1", "1", "3.13%", "NO"),
    @("4324", "1", "3", "2", "16.96%", "NO"),
    @("4332", "", "", "", "", "NO"),
    @("4338", "1", "6", "0", "17.91%", "NO"),
    @("4342", "1", "0", "0", "", "NO"),
    @("4345", "1", "6", "1", "24.52%", "NO"),
    @("4350", "1", "11", "1", "41.89%", "NO"),
    @("4353", "", "", "", "", "NO"),
    @("4385", "", "", "", "", "NO"),
    @("4387", "2", "8", "3", "26.36%", "NO"),
    @("4388", "2", "8", "1", "24.37%", "NO"),
    @("4398", "3", "4", "0", "18.43%", "NO"),
    @("4399", "5", "6", "3", "23.53%", "YES"),
    @("4400", "2", "2", "0", "6.57%", "NO"),
    @("4402", "", "", "", "", "NO"),
    @("4406", "", "", "", "", "NO"),
    @("4410", "5", "9", "2", "37.84%", "NO"),
    @("4435", "", "", "", "", "NO"),
    @("4436", "", "", "", "", "NO")
)

for ($i = 0; $i -lt $newData.Length; $i++) {
    $row = $newData[$i]
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne "") {
        $ws.Cells.Item($r, 2).Value = [double]$row[1]
    }
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
